# Generate Report for Handoff
# Updates the localization-status report after a new handoff generation run:
#  - Priority for the still-pending (non "Handed back") rows flips from "low" to "ht"
#  - The Latest Handoff Datetime for the zh-cn locale reflects the new handoff run
#  - The Latest HO Xliff Generate Date summary (Overview sheet, also mirrored onto
#    the de-de sheet's Latest Handoff Datetime column) advances accordingly

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

for ($row = 4; $row -le 7; $row++) {
    $overview.Range("G$row").Value = "2016-08-27 16:31:28"

    $zhcn.Range("E$row").Value = "ht"
    $zhcn.Range("H$row").Value = "2016-08-27 16:31:23"

    $dede.Range("E$row").Value = "ht"
    $dede.Range("H$row").Value = "2016-08-27 16:31:28"
}
